$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths / styles (cols B widened, A/C/D get the "style=2" text
#    number-format column style; col B loses its custom "style=3")
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 23.7109375

# ---------------------------------------------------------------------------
# 2. Re-style the data cells.
#    Column A and C (AREA_CODE / BRANCH_CODE) need the plain bordered cell
#    with a Text ("@") number format.
#    Column B, D, E, F (AREA_NAME / BRANCH_NAME / OS_AMOUNT / WO_AMOUNT) need
#    the plain bordered cell with the General number format (this drops the
#    special bold/shaded font+fill that used to live on column D and on the
#    amount columns).
#    We copy formats from A2 (already "plain + border") via PasteSpecial so
#    we reuse the existing style entries instead of inventing new ones, then
#    layer on the "@" text format where required.
# ---------------------------------------------------------------------------
$plain = $ws.Range("A2")
$plain.Copy()
$ws.Range("B2:B4,D2:D4,E2:E4,F2:F4").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:A4,C2:C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("C2:C4").NumberFormat = "@"
$ws.Range("B2:B4,D2:D4,E2:E4,F2:F4").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 3. New cell values (the "real data")
# ---------------------------------------------------------------------------
# Row 2 - Ashkona Branch
$ws.Range("A2").Value = "01"
$ws.Range("B2").Value = "Dhaka Area"
$ws.Range("C2").Value = 124
$ws.Range("D2").Value = "Ashkona Branch"
$ws.Range("E2").Value = 2000
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 44196

# Row 3 - Chandgaon Branch
$ws.Range("A3").Value = "03"
$ws.Range("B3").Value = "Chattogram Area 1"
$ws.Range("C3").Value = "007"
$ws.Range("D3").Value = "Chandgaon Branch"
$ws.Range("E3").Value = 1500
$ws.Range("F3").Value = 1500
$ws.Range("G3").Value = 44196

# Row 4 - Feni SME Branch
$ws.Range("A4").Value = "03"
$ws.Range("B4").Value = "Chattogram Area 1"
$ws.Range("C4").Value = "012"
$ws.Range("D4").Value = "Feni SME Branch"
$ws.Range("E4").Value = 1200
$ws.Range("F4").Value = 2000
$ws.Range("G4").Value = 44196

# ---------------------------------------------------------------------------
# 4. Selection marker left behind in the sheet view by the author
# ---------------------------------------------------------------------------
$ws.Range("L25").Select()
